$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook tracks a batch of "downtime" events identified by a uuid in
# column G. This edit rolls the existing rows (2-25) onto a new uuid and
# appends a fresh copy of the first block of rows (2-9) underneath (rows
# 26-33) stamped with that same new uuid - i.e. "version 2" of the data.

$newGuid = "d85554b9-776c-49d1-bdf2-3016191cd60b"

# 1) Re-stamp every existing data row (2-25) with the new uuid.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 7).Value = $newGuid
}

# 2) Append rows 26-33: a duplicate of rows 2-9 (which now already carry the
#    new uuid), preserving values and number formats.
$ws.Range("A2:H9").Copy($ws.Range("A26"))
